$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '319.93'
Set-TextValue $ws.Range('E2') '3.68%'
Set-TextValue $ws.Range('D3') '41.32'
Set-TextValue $ws.Range('E3') '1.51%'
Set-TextValue $ws.Range('D4') '5.244'
Set-TextValue $ws.Range('E4') '2.22%'
Set-TextValue $ws.Range('D5') '0.07746'
Set-TextValue $ws.Range('E5') '1.78%'
Set-TextValue $ws.Range('B6') 'FTXToken'
Set-TextValue $ws.Range('C6') 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D6') '1.692'
Set-TextValue $ws.Range('E6') '5.11%'
Set-TextValue $ws.Range('B7') 'MXToken'
Set-TextValue $ws.Range('C7') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D7') '0.9442'
Set-TextValue $ws.Range('E7') '3.97%'
Set-TextValue $ws.Range('B8') 'BTSEToken'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D8') '2.425'
Set-TextValue $ws.Range('E8') '-1.61%'
Set-TextValue $ws.Range('B9') 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D9') '0.1244'
Set-TextValue $ws.Range('E9') '-3.16%'
Set-TextValue $ws.Range('B10') 'WazirX'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D10') '0.1838'
Set-TextValue $ws.Range('E10') '1.58%'
Set-TextValue $ws.Range('B11') 'MandalaExchangeToken'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D11') '0.09221'
Set-TextValue $ws.Range('E11') '1.73%'
Set-TextValue $ws.Range('B12') 'BitrueCoin'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D12') '0.04325'
Set-TextValue $ws.Range('E12') '0.10%'
Set-TextValue $ws.Range('B13') 'BitMartToken'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D13') '0.1050'
Set-TextValue $ws.Range('E13') '0.52%'
Set-TextValue $ws.Range('B14') 'BitForexToken'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D14') '0.001289'
Set-TextValue $ws.Range('E14') '2.84%'
Set-TextValue $ws.Range('B15') 'TigerCash'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D15') '0.005957'
Set-TextValue $ws.Range('E15') '1.15%'
Set-TextValue $ws.Range('B16') 'UpBots'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws.Range('D16') '0.007491'
Set-TextValue $ws.Range('E16') '1,897.31%'
Set-TextValue $ws.Range('B17') 'HotbitToken'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D17') '0.004121'
Set-TextValue $ws.Range('E17') '2.20%'
Set-TextValue $ws.Range('E18') '-0.11%'
Set-TextValue $ws.Range('B19') 'GateToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D19') '4.340'
Set-TextValue $ws.Range('E19') '1.36%'
Set-TextValue $ws.Range('B20') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D20') '0.3329'
Set-TextValue $ws.Range('E20') '0.46%'
Set-TextValue $ws.Range('B21') 'MCDex'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D21') '7.770'
Set-TextValue $ws.Range('E21') '11.19%'
Set-TextValue $ws.Range('B22') 'ProBitToken'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D22') '0.1354'
Set-TextValue $ws.Range('E22') '-3.01%'
Set-TextValue $ws.Range('B23') 'ZBToken'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D23') '0.2824'
Set-TextValue $ws.Range('E23') '4.37%'
Set-TextValue $ws.Range('B24') 'CoinExToken'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D24') '0.04036'
Set-TextValue $ws.Range('E24') '-0.10%'
Set-TextValue $ws.Range('B25') 'BitKan'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D25') '0.001265'
Set-TextValue $ws.Range('E25') '-0.72%'
Set-TextValue $ws.Range('E26') '-0.06%'
Set-TextValue $ws.Range('D38') '0.02554'
Set-TextValue $ws.Range('E38') '5.48%'
Set-TextValue $ws.Range('D39') '0.05354'
Set-TextValue $ws.Range('E39') '2.64%'
Set-TextValue $ws.Range('D40') '0.007757'
Set-TextValue $ws.Range('E40') '-1.02%'
Set-TextValue $ws.Range('D41') '0.1317'
Set-TextValue $ws.Range('E41') '1.47%'
Set-TextValue $ws.Range('D42') '0.007360'
Set-TextValue $ws.Range('E42') '8.11%'
Set-TextValue $ws.Range('D43') '0.001992'
Set-TextValue $ws.Range('E43') '4.88%'
Set-TextValue $ws.Range('D44') '0.008350'
Set-TextValue $ws.Range('E44') '13.01%'
Set-TextValue $ws.Range('D45') '0.3178'
Set-TextValue $ws.Range('E45') '-5.09%'
Set-TextValue $ws.Range('D46') '0.00006718'
Set-TextValue $ws.Range('E46') '-2.80%'
Set-TextValue $ws.Range('E47') '-0.09%'
Set-TextValue $ws.Range('D48') '0.1997'
Set-TextValue $ws.Range('E48') '82.91%'
Set-TextValue $ws.Range('D49') '0.004203'
Set-TextValue $ws.Range('E49') '39.99%'
Set-TextValue $ws.Range('E50') '-0.09%'
Set-TextValue $ws.Range('E51') '-0.09%'

Write-Host "Applied $($ws.Range('A1').Worksheet.Name) updates"
